# Fix to Issue 23490 : Middle name field not included in HDRL test request template
#
# Insert a new "Middle Name" column between "First Name" (col C) and
# "Birth Date" (col D), shifting every later column one place to the
# right (Birth Date, SSN, Draw Date, FMP, DUC, SOT, DOD Id).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D - this shifts D..J (Birth Date..DOD Id) to E..K
# and copies formatting (style) from the neighbouring column, same as
# Excel's UI "Insert" on a column header.
$ws.Columns("D").Insert()

# Populate the new header cell.
$ws.Cells.Item(1, 4).Value = "Middle Name"

# Match the column's on-disk width to the new header's contents.
$ws.Columns("D").ColumnWidth = 10.83

# Leave the same cell selected/active as in the saved workbook.
$ws.Range("E3").Select() | Out-Null
